{"js": "// Office.js (Word JavaScript API) script.\n// Applies the edit described by the diff:\n//  1. Remove the \"Meta description: ...\" paragraph that sits right after the\n//     title heading at the top of the document.\n//  2. Insert a new bold paragraph \"Play Fat Santa Slot Free - Review &\n//     Ratings 2021\" right before the final (italic) paragraph.\n//  3. Replace the text of that final paragraph (previously an image-prompt\n//     sentence) with the meta-description sentence, keeping its italic run\n//     formatting intact.\n\nconst body = context.document.body;\n\n// --- Step 1: delete the \"Meta description\" paragraph -----------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet metaPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Meta description\") !== -1) {\n    metaPara = paragraphs.items[i];\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// --- Step 2: insert the new bold paragraph before the last paragraph ---\n// Re-load the paragraph collection since the body changed above.\nconst paragraphsAfterDelete = body.paragraphs;\nparagraphsAfterDelete.load(\"items/text\");\nawait context.sync();\n\nconst itemsAfterDelete = paragraphsAfterDelete.items;\nconst lastParaBeforeInsert = itemsAfterDelete[itemsAfterDelete.length - 1];\nconst secondToLastPara = itemsAfterDelete[itemsAfterDelete.length - 2];\n\n// Insert a brand-new bold paragraph right after the second-to-last\n// paragraph (i.e. immediately before the last paragraph), using raw OOXML\n// so we get the exact run layout (leading empty run + bold run) used\n// elsewhere in this document.\nconst insertionRange = secondToLastPara.getRange(\"End\");\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +\n  '<w:t>Play Fat Santa Slot Free - Review &amp; Ratings 2021</w:t>' +\n  '</w:r></w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ninsertionRange.insertOoxml(newParagraphOoxml, \"After\");\nawait context.sync();\n\n// --- Step 3: replace the text of the (new) final paragraph -------------\n// Reload once more so we grab the paragraph that is now truly last (the\n// original \"Create a feature image...\" paragraph), not the one we just\n// inserted.\nconst paragraphsFinal = body.paragraphs;\nparagraphsFinal.load(\"items/text\");\nawait context.sync();\n\nconst finalItems = paragraphsFinal.items;\nconst trueLastPara = finalItems[finalItems.length - 1];\ntrueLastPara.insertText(\n  \"Read our review of Fat Santa slot game. Discover the game features and play it for free to win real money. Find out if it's worth playing.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the edit described by the diff:\n#  1. Remove the \"Meta description: ...\" paragraph that sits right after the\n#     title heading at the top of the document.\n#  2. Insert a new bold paragraph \"Play Fat Santa Slot Free - Review &\n#     Ratings 2021\" right before the final (italic) paragraph.\n#  3. Replace the text of that final paragraph (previously an image-prompt\n#     sentence) with the meta-description sentence, keeping its italic run\n#     formatting intact.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Meta description\" paragraph -----------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Meta description*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- Step 2: insert the new bold paragraph before the last paragraph ---\n$count = $d.Paragraphs.Count\n$secondToLast = $d.Paragraphs.Item($count - 1)\n# Collapse to a point just inside the end of the second-to-last paragraph\n# (one character before its paragraph mark). Inserting raw OOXML exactly at\n# a paragraph boundary can consume/replace the neighbouring paragraph in\n# this engine, so we target a point a single character earlier, which\n# cleanly splits the paragraph mark off into a new paragraph instead.\n$pos = $secondToLast.Range.End - 1\n$insertionPoint = $d.Range($pos, $pos)\n$newParagraphXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fat Santa Slot Free - Review &amp; Ratings 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n[void]$insertionPoint.InsertXML($newParagraphXml)\n\n# --- Step 3: replace the text of the (now true) final paragraph --------\n$count2 = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count2)\n$textRange = $lastPara.Range\n# Exclude the trailing paragraph mark so the replacement keeps living in the\n# same paragraph (and keeps its existing italic run formatting).\n[void]$textRange.MoveEnd(1, -1)\n$textRange.Text = \"Read our review of Fat Santa slot game. Discover the game features and play it for free to win real money. Find out if it's worth playing.\"\n"}
